$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 6.8499999899999997
$ws.Range("C2").Value = 7.2999999899999999
$ws.Range("D2").Value = 6.8499999899999997
$ws.Range("E2").Value = 7.2999999899999999

$ws.Range("B3").Value = 6.8499999899999997
$ws.Range("C3").Value = 7.2999999899999999
$ws.Range("D3").Value = 6.8499999899999997
$ws.Range("E3").Value = 7.2999999899999999

$ws.Range("B1:E3").Select()
